$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @("IMX-USD", "TAO-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 439
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
